# Generate Report for Handback
# Update the timestamp values for the "a2c1d7f3-05cc-4888-95b8-cccdb897c336" row
# across the Overview, zh-cn, and de-de worksheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-28 02:44:59"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-28 02:44:54"
$wsZhCn.Range("K4").Value = "2016-08-28 02:45:28"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-08-28 02:45:35"
